$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = "19.9648155496385"
    3 = "16.336820679251357"
    4 = "19.109729235928555"
    5 = "15.36459887606487"
    6 = "14.481527814310157"
    7 = "17.20350490714225"
    8 = "21.623741553475075"
}

foreach ($row in $values.Keys) {
    $text = "ReturnTuple(sdnn=" + $values[$row] + ")"
    foreach ($col in @("D", "E", "F", "G")) {
        $ws.Range($col + $row).Value = $text
    }
}
